$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.433.57"
Set-TextValue "E2" "  +0.37%  "
Set-TextValue "D3" "2.284.92"
Set-TextValue "E3" "  -0.97%  "
Set-TextValue "E4" "  -0.16%  "
Set-TextValue "D5" "324.70"
Set-TextValue "E5" "  +2.27%  "
Set-TextValue "D6" "102.85"
Set-TextValue "E6" "  -2.15%  "
Set-TextValue "D7" "0.622"
Set-TextValue "E7" "  -1.10%  "
Set-TextValue "E8" "  -0.04%  "
Set-TextValue "D9" "0.608"
Set-TextValue "D10" "39.91"
Set-TextValue "E10" "  +0.06%  "
Set-TextValue "D11" "0.0907"
Set-TextValue "E11" "  -0.41%  "
Set-TextValue "D12" "8.33"
Set-TextValue "E12" "  -1.44%  "
Set-TextValue "D13" "0.107"
Set-TextValue "E13" "  -0.27%  "
Set-TextValue "D14" "0.970"
Set-TextValue "E14" "  -0.89%  "
Set-TextValue "D15" "15.09"
Set-TextValue "E15" "  -2.64%  "
Set-TextValue "D16" "2.628.97"
Set-TextValue "E16" "  -1.00%  "
Set-TextValue "D17" "2.285.74"
Set-TextValue "E17" "  -1.09%  "
Set-TextValue "D18" "42.249.51"
Set-TextValue "E18" "  +0.07%  "
Set-TextValue "D19" "7.37"
Set-TextValue "E19" "  -5.34%  "
Set-TextValue "E20" "  -0.67%  "
Set-TextValue "D21" "12.95"
Set-TextValue "E21" "  +29.07%  "
Set-TextValue "D22" "3.64"
Set-TextValue "E22" "  +2.28%  "
Set-TextValue "D23" "73.02"
Set-TextValue "E23" "  -1.21%  "
Set-TextValue "D24" "267.63"
Set-TextValue "E24" "  -6.42%  "
Set-TextValue "E26" "  -0.21%  "
Set-TextValue "D27" "10.86"
Set-TextValue "E27" "  -1.26%  "
Set-TextValue "E28" "  +4.22%  "
Set-TextValue "B29" "InjectiveProtocol"
Set-TextValue "C29" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D29" "38.24"
Set-TextValue "E29" "  +7.33%  "
Set-TextValue "B30" "EthereumClassic"
Set-TextValue "C30" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D30" "22.41"
Set-TextValue "E30" "  -4.35%  "
Set-TextValue "D31" "164.38"
Set-TextValue "E31" "  -0.52%  "
Set-TextValue "D32" "6.14"
Set-TextValue "E32" "  +3.37%  "
Set-TextValue "D33" "0.0881"
Set-TextValue "E33" "  -0.40%  "
Set-TextValue "E34" "  +0.45%  "
Set-TextValue "D35" "2.52"
Set-TextValue "E35" "  -13.75%  "
Set-TextValue "D36" "0.113"
Set-TextValue "E36" "  -4.74%  "
Set-TextValue "E37" "  -1.67%  "
Set-TextValue "E38" "  -0.19%  "
Set-TextValue "D39" "3.72"
Set-TextValue "E39" "  +1.84%  "
Set-TextValue "E40" "  -7.02%  "
Set-TextValue "D41" "1.52"
Set-TextValue "E41" "  +1.33%  "
Set-TextValue "D42" "69.57"
Set-TextValue "E42" "  -2.34%  "
Set-TextValue "E43" "  -0.35%  "
Set-TextValue "E44" "  -1.07%  "
Set-TextValue "B45" "BitcoinSV"
Set-TextValue "C45" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D45" "91.20"
Set-TextValue "E45" "  -10.56%  "
Set-TextValue "B46" "Celestia"
Set-TextValue "C46" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D46" "12.35"
Set-TextValue "E46" "  +1.35%  "
Set-TextValue "D47" "113.12"
Set-TextValue "E47" "  -3.17%  "
Set-TextValue "D48" "80.00"
Set-TextValue "E48" "  +1.41%  "
Set-TextValue "D49" "8.93"
Set-TextValue "E49" "  -2.83%  "
Set-TextValue "D50" "5.21"
Set-TextValue "E50" "  -2.77%  "
Set-TextValue "D51" "1.588.84"
Set-TextValue "E51" "  +1.68%  "
